$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 4.3
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.21
$ws.Range("R2").Value = 1.55
$ws.Range("S2").Value = 2.52
$ws.Range("T2").Value = 1.69
$ws.Range("U2").Value = 2.18
$ws.Range("V2").Value = 1.2
$ws.Range("X2").Value = 24
$ws.Range("Y2").Value = 25
$ws.Range("Z2").Value = 980
$ws.Range("AA2").Value = 150
$ws.Range("AB2").Value = 12
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 80
$ws.Range("AF2").Value = 12.5
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 75
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 19.5
$ws.Range("AL2").Value = 30
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 7.8
$ws.Range("AO2").Value = 75

$ws.Range("J3").Value = 3.95
$ws.Range("Q3").Value = 1.91

$ws.Range("F5").Value = 1.41
$ws.Range("I5").Value = 8.800000000000001
$ws.Range("K5").Value = 6.6
$ws.Range("P5").Value = 3.45

$ws.Range("G6").Value = 2.2
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 4.3
$ws.Range("J6").Value = 3.9

$ws.Range("F7").Value = 1.44
$ws.Range("G7").Value = 1.61
$ws.Range("H7").Value = 6.4
$ws.Range("I7").Value = 8.800000000000001
$ws.Range("K7").Value = 5.9
$ws.Range("P7").Value = 2.9
$ws.Range("Q7").Value = 1.43

$ws.Range("I8").Value = 2.8
$ws.Range("P8").Value = 1.97
$ws.Range("Q8").Value = 1.88
